$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.398.07'
$ws.Range('D3').Value = '1.822.52'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'313.90"
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  +1.91%  '
$ws.Range('D8').Value = "'0.3750"
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('D9').Value = "'0.07501"
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').Value = "'0.8862"
$ws.Range('E10').Value = '  +4.94%  '
$ws.Range('D11').Value = "'21.10"
$ws.Range('E11').Value = '  +1.90%  '
$ws.Range('D12').Value = '1.824.03'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = "'6.755"
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').Value = "'93.95"
$ws.Range('E14').Value = '  +4.49%  '
$ws.Range('D15').Value = "'5.413"
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = "'0.07103"
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = "'1.001"
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = "'0.000008782"
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = "'15.17"
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').Value = '27.398.18'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').Value = "'5.310"
$ws.Range('E22').Value = '  +3.10%  '
$ws.Range('D23').Value = "'10.93"
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').Value = '2.056.40'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('D26').Value = "'2.363"
$ws.Range('E26').Value = '  +7.20%  '
$ws.Range('D27').Value = "'151.59"
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'18.57"
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('D29').Value = "'5.364"
$ws.Range('E29').Value = '  +2.48%  '
$ws.Range('D30').Value = "'118.19"
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('D31').Value = "'0.08848"
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').Value = "'0.7871"
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('D33').Value = "'1.201"
$ws.Range('E33').Value = '  +1.69%  '
$ws.Range('D34').Value = "'4.535"
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').Value = "'2.918"
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('D36').Value = "'0.9998"
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = "'1.111"
$ws.Range('E37').Value = '  +1.37%  '
$ws.Range('E38').Value = '  +2.61%  '
$ws.Range('D39').Value = "'0.05325"
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('D40').Value = "'7.375"
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('D41').Value = "'0.5322"
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('D42').Value = "'0.1722"
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('D43').Value = "'2.856"
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('D44').Value = "'2.307"
$ws.Range('E44').Value = '  +19.36%  '
$ws.Range('D45').Value = "'8.724"
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('D46').Value = "'0.5094"
$ws.Range('E46').Value = '  +5.77%  '
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = "'105.70"
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'1.700"
$ws.Range('E49').Value = '  +2.47%  '
$ws.Range('D50').Value = "'0.9999"
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('E51').Value = '  +0.78%  '
